$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook now includes a new "ECs" target cluster, which adds one new
# row for each existing "Sending cluster" (FAPs, sCs) -- growing the table
# from 4 to 6 data rows. The whole LR-pair table was regenerated with the
# extra cluster included, so every numeric value changes too, not just the
# new rows. Rewrite every data cell (rows 2-7) directly with its final
# value rather than inserting rows (Rows.Insert() here would copy the bold
# header's style down onto the new row).

$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(2,2).Value2 = "Fgf7"
$ws.Cells.Item(2,3).Value2 = "Fgfr3"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 13.79659733333333
$ws.Cells.Item(2,8).Value2 = 41.389792
$ws.Cells.Item(2,9).Value2 = 0.9485830781324925
$ws.Cells.Item(2,10).Value2 = 0.9485830781324925
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 4.220261333333333
$ws.Cells.Item(2,14).Value2 = 12.660784
$ws.Cells.Item(2,15).Value2 = 0.6739259863235564
$ws.Cells.Item(2,16).Value2 = 0.6739259863235564
$ws.Cells.Item(2,17).Value2 = 58.22524625743644
$ws.Cells.Item(2,18).Value2 = 524.0272163169279
$ws.Cells.Item(2,19).Value2 = 0.6392747865402751
$ws.Cells.Item(2,20).Value2 = 0.6392747865402751
$ws.Cells.Item(3,1).Value2 = "FAPs"
$ws.Cells.Item(3,2).Value2 = "Fgf7"
$ws.Cells.Item(3,3).Value2 = "Fgfr3"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 3
$ws.Cells.Item(3,6).Value2 = 1
$ws.Cells.Item(3,7).Value2 = 13.79659733333333
$ws.Cells.Item(3,8).Value2 = 41.389792
$ws.Cells.Item(3,9).Value2 = 0.9485830781324925
$ws.Cells.Item(3,10).Value2 = 0.9485830781324925
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 1.296447666666667
$ws.Cells.Item(3,14).Value2 = 3.889343
$ws.Cells.Item(3,15).Value2 = 0.2070274097896007
$ws.Cells.Item(3,16).Value2 = 0.2070274097896007
$ws.Cells.Item(3,17).Value2 = 17.88656642073956
$ws.Cells.Item(3,18).Value2 = 160.979097786656
$ws.Cells.Item(3,19).Value2 = 0.1963826976360163
$ws.Cells.Item(3,20).Value2 = 0.1963826976360163
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Fgf7"
$ws.Cells.Item(4,3).Value2 = "Fgfr3"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 13.79659733333333
$ws.Cells.Item(4,8).Value2 = 41.389792
$ws.Cells.Item(4,9).Value2 = 0.9485830781324925
$ws.Cells.Item(4,10).Value2 = 0.9485830781324925
$ws.Cells.Item(4,11).Value2 = 3
$ws.Cells.Item(4,12).Value2 = 1
$ws.Cells.Item(4,13).Value2 = 0.745494
$ws.Cells.Item(4,14).Value2 = 2.236482
$ws.Cells.Item(4,15).Value2 = 0.119046603886843
$ws.Cells.Item(4,16).Value2 = 0.119046603886843
$ws.Cells.Item(4,17).Value2 = 10.285280532416
$ws.Cells.Item(4,18).Value2 = 92.56752479174401
$ws.Cells.Item(4,19).Value2 = 0.1129255939562011
$ws.Cells.Item(4,20).Value2 = 0.1129255939562011
$ws.Cells.Item(5,1).Value2 = "sCs"
$ws.Cells.Item(5,2).Value2 = "Fgf7"
$ws.Cells.Item(5,3).Value2 = "Fgfr3"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 0.7478296666666666
$ws.Cells.Item(5,8).Value2 = 2.243489
$ws.Cells.Item(5,9).Value2 = 0.05141692186750751
$ws.Cells.Item(5,10).Value2 = 0.05141692186750751
$ws.Cells.Item(5,11).Value2 = 2
$ws.Cells.Item(5,12).Value2 = 0.6666666666666666
$ws.Cells.Item(5,13).Value2 = 4.220261333333333
$ws.Cells.Item(5,14).Value2 = 12.660784
$ws.Cells.Item(5,15).Value2 = 0.6739259863235564
$ws.Cells.Item(5,16).Value2 = 0.6739259863235564
$ws.Cells.Item(5,17).Value2 = 3.156036626152888
$ws.Cells.Item(5,18).Value2 = 28.404329635376
$ws.Cells.Item(5,19).Value2 = 0.03465119978328124
$ws.Cells.Item(5,20).Value2 = 0.03465119978328124
$ws.Cells.Item(6,1).Value2 = "sCs"
$ws.Cells.Item(6,2).Value2 = "Fgf7"
$ws.Cells.Item(6,3).Value2 = "Fgfr3"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 3
$ws.Cells.Item(6,6).Value2 = 1
$ws.Cells.Item(6,7).Value2 = 0.7478296666666666
$ws.Cells.Item(6,8).Value2 = 2.243489
$ws.Cells.Item(6,9).Value2 = 0.05141692186750751
$ws.Cells.Item(6,10).Value2 = 0.05141692186750751
$ws.Cells.Item(6,11).Value2 = 3
$ws.Cells.Item(6,12).Value2 = 1
$ws.Cells.Item(6,13).Value2 = 1.296447666666667
$ws.Cells.Item(6,14).Value2 = 3.889343
$ws.Cells.Item(6,15).Value2 = 0.2070274097896007
$ws.Cells.Item(6,16).Value2 = 0.2070274097896007
$ws.Cells.Item(6,17).Value2 = 0.969522026414111
$ws.Cells.Item(6,18).Value2 = 8.725698237727
$ws.Cells.Item(6,19).Value2 = 0.01064471215358436
$ws.Cells.Item(6,20).Value2 = 0.01064471215358436
$ws.Cells.Item(7,1).Value2 = "sCs"
$ws.Cells.Item(7,2).Value2 = "Fgf7"
$ws.Cells.Item(7,3).Value2 = "Fgfr3"
$ws.Cells.Item(7,4).Value2 = "sCs"
$ws.Cells.Item(7,5).Value2 = 3
$ws.Cells.Item(7,6).Value2 = 1
$ws.Cells.Item(7,7).Value2 = 0.7478296666666666
$ws.Cells.Item(7,8).Value2 = 2.243489
$ws.Cells.Item(7,9).Value2 = 0.05141692186750751
$ws.Cells.Item(7,10).Value2 = 0.05141692186750751
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 0.745494
$ws.Cells.Item(7,14).Value2 = 2.236482
$ws.Cells.Item(7,15).Value2 = 0.119046603886843
$ws.Cells.Item(7,16).Value2 = 0.119046603886843
$ws.Cells.Item(7,17).Value2 = 0.557502529522
$ws.Cells.Item(7,18).Value2 = 5.017522765698
$ws.Cells.Item(7,19).Value2 = 0.006121009930641925
$ws.Cells.Item(7,20).Value2 = 0.006121009930641925
